$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Fitness (column C) values for rows 2 through 180 (179 values)
$values = @(
13819,12139,11158,11158,11158,10370,10301,8779,8779,8779,8779,8779,8779,8779,8744,8744,8497,8497,8344,8344,8344,8344,8344,8344,8344,8344,8344,8344,8344,8344,8344,8344,8344,8344,8344,8344,8344,8344,8344,7899,7899,7899,7899,7899,7899,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7807,7618,7618,7343,7343,7343,7343,7343,7343,7343,7343,7343,7343,7343,7343,7343,7343,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293
)

$arr = New-Object "object[,]" $values.Length,1
for ($i = 0; $i -lt $values.Length; $i++) {
    $arr[$i, 0] = $values[$i]
}

$ws.Range("C2:C180").Value = $arr
